$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6919785141944885
$ws.Range("B1").Value = 1.046266078948975
$ws.Range("C1").Value = 3.930750608444214
$ws.Range("D1").Value = 3.549507856369019
$ws.Range("E1").Value = 1.977318406105042
